$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.357.94'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '3.140.20'
$ws.Range("E3").Value = '  -1.09%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.56'
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.95'
$ws.Range("E6").Value = '  -2.80%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").Value = '  -4.73%  '

$ws.Range("D9").Value = '3.153.56'
$ws.Range("E9").Value = '  -1.01%  '

$ws.Range("E10").Value = '  -3.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.60'
$ws.Range("E11").Value = '  -3.15%  '

$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("D13").Value = '3.684.73'
$ws.Range("E13").Value = '  -1.50%  '

$ws.Range("E14").Value = '  -1.44%  '

$ws.Range("D15").Value = '64.387.72'
$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.10'
$ws.Range("E16").Value = '  -1.11%  '

$ws.Range("D17").Value = '3.140.36'
$ws.Range("E17").Value = '  -1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000155'
$ws.Range("E18").Value = '  -2.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '402.63'
$ws.Range("E19").Value = '  -3.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.26'
$ws.Range("E20").Value = '  -1.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.57'
$ws.Range("E21").Value = '  -2.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.08'
$ws.Range("E22").Value = '  -0.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.84'
$ws.Range("E24").Value = '  -2.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.486'
$ws.Range("E25").Value = '  -0.53%  '

$ws.Range("E26").Value = '  -4.47%  '

$ws.Range("E27").Value = '  -3.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.82'
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.09%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.81'
$ws.Range("E31").Value = '  -1.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.26'
$ws.Range("E32").Value = '  -2.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '161.19'
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.88'
$ws.Range("E34").Value = '  -4.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.28'
$ws.Range("E35").Value = '  -0.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.12'
$ws.Range("E36").Value = '  -1.95%  '

$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("E38").Value = '  -1.56%  '

$ws.Range("D39").Value = '2.638.47'
$ws.Range("E39").Value = '  -3.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.74'
$ws.Range("E40").Value = '  -2.93%  '

$ws.Range("E41").Value = '  -3.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.53'
$ws.Range("E42").Value = '  -1.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.690'
$ws.Range("E43").Value = '  -3.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0615'
$ws.Range("E44").Value = '  -1.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.43'
$ws.Range("E45").Value = '  -3.20%  '

$ws.Range("E46").Value = '  -3.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.20'
$ws.Range("E47").Value = '  -2.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '287.59'
$ws.Range("E48").Value = '  -1.94%  '

$ws.Range("E49").Value = '  -0.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0979'
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.91'
$ws.Range("E51").Value = '  -4.74%  '
